$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1Values = @(
    "-7.2492138255929334E-3", "-8.1406616791069343E-3", "-1.0476259361751117E-2", "-2.556858610363141E-2", "-1.0125570160998762E-2", "-1.3172471437025322E-2", "-1.2641673179481956E-2", "-1.1937456888638792E-2", "-1.378804864623317E-2", "-7.247545319591083E-3",
    "-1.646151114502167E-2", "-1.3512618722307355E-2", "-1.251861593687624E-2", "-5.571747480154661E-3", "-1.4462009470536831E-2", "-1.6551990704047538E-2", "-5.8200881485237161E-3", "-1.3496884892132317E-2", "-1.3872085034304597E-2", "-5.5951081168035188E-3",
    "-8.8237923253074221E-3", "-5.8510575717777916E-3", "-1.5427601051563603E-2", "-8.7213563534530598E-3", "-1.1082330541946834E-2", "-9.9556622633120966E-3", "-6.2275913804645828E-3", "-1.3861635660099954E-2", "-9.3971134133796173E-3", "-1.1075932117370757E-2",
    "-6.6305612759871065E-3", "-6.893430232916253E-3", "-1.3389139047476841E-2", "-5.8670421861891729E-3", "-1.2906145175809493E-2", "-1.2169101628959917E-2", "-1.4189892947085466E-2", "-1.093759161863948E-2", "-9.7296514478636895E-3", "-1.1266775644244725E-2",
    "-8.8178120656835571E-3", "-1.5547878187930932E-2", "-1.1552749219122228E-2", "-1.3598249296308396E-2", "-7.8437667642468172E-3", "-7.0830787390293972E-3", "-7.8593407876522085E-3", "-1.2574802918090995E-2", "-1.6249698049970198E-2", "-7.5256715133203709E-3",
    "-1.481900985676071E-2", "-1.1592402383972571E-2", "-1.1024348785214182E-2", "-9.7984816463209702E-3", "-2.0626752990217063E-2", "-1.7219725455981397E-2", "-8.0782716237869957E-3", "-3.8681819646098209E-3", "-8.2957899077579066E-3", "-1.425184922553157E-2",
    "-1.3169274263205433E-2", "-1.5260177747871475E-2", "-8.7700092893821004E-3", "-1.2856714945514405E-2", "-1.205544726356563E-2", "-9.591909969939352E-3", "-9.8161960700519565E-3", "-9.5041388609923214E-3", "-9.9404947335626718E-3", "-1.2550208508715429E-2",
    "-7.0414880134418087E-3", "-1.172784490709276E-2", "-6.7806898315991657E-3", "-9.6182641589219482E-3", "-9.3171513970979301E-3", "-1.2145755516404714E-2", "-1.1174883837687579E-2", "-1.4894073434866317E-2", "-2.0639696556917504E-2", "-8.7641346547216951E-3",
    "-9.8057601929620666E-3", "-2.4627395122719914E-2", "-9.1017642920785506E-3", "-1.1337532676992989E-2", "-9.4551228883978603E-3", "-1.3723318436480536E-2", "-1.9702816905792078E-2", "-1.7959128558727371E-2", "-1.289130835119561E-2", "-1.0958252114243094E-2",
    "-1.5916930086667076E-2", "-7.3190747732922948E-3", "-1.6179571103717798E-2", "-8.6265468699684138E-3", "-9.2314960547410924E-3", "-1.2026689574848769E-2", "-7.8644794802026196E-3", "-1.4603943371106199E-2", "-1.8644272319801463E-2", "-1.5127456752527305E-2"
)

$row2Values = @(
    "-4.8350637807840205E-3", "-5.5035107910001652E-3", "-7.0806624105755438E-3", "-1.7003156279647283E-2", "-6.7238527716587705E-3", "-8.7383570410764638E-3", "-8.4962159653004087E-3", "-7.9760443121369533E-3", "-9.0628326991734447E-3", "-4.8656856824941235E-3",
    "-1.0905538019935926E-2", "-8.9766993561141307E-3", "-8.4943799651096257E-3", "-3.822264160006238E-3", "-9.5428738897912809E-3", "-1.1080266826027323E-2", "-3.9794242965338856E-3", "-9.0981358140962507E-3", "-9.1325682290481509E-3", "-3.7450789119077184E-3",
    "-5.8306731614316846E-3", "-3.8748309883264997E-3", "-1.0186842602661658E-2", "-5.8524473823311721E-3", "-7.4389505303242738E-3", "-6.7097231218909755E-3", "-4.0054858365030496E-3", "-9.3832609007575632E-3", "-6.1680565897098763E-3", "-7.5214073125262565E-3",
    "-4.3194508629079862E-3", "-4.7068760719380541E-3", "-8.8840634940467126E-3", "-3.8239407701375684E-3", "-8.5763349988333751E-3", "-8.0521151536310509E-3", "-9.5984490206781903E-3", "-7.2975704171883505E-3", "-6.4506100043109508E-3", "-7.4436883520902176E-3",
    "-6.0355367234764534E-3", "-1.0494989162986046E-2", "-7.8348364910579039E-3", "-9.0779372236278171E-3", "-5.2468434065360515E-3", "-4.6596777287269417E-3", "-5.2181775809046057E-3", "-8.3367713879752142E-3", "-1.0765197446281318E-2", "-4.9456614482125301E-3",
    "-9.6641012471754256E-3", "-7.7837635561044099E-3", "-7.3714654808106515E-3", "-6.4278460062444209E-3", "-1.3839566425071307E-2", "-1.1514527958063559E-2", "-5.3309719625264391E-3", "-2.6389272263768458E-3", "-5.4140324822457035E-3", "-9.6406738261237274E-3",
    "-8.8322995177219973E-3", "-1.036544253982016E-2", "-5.864174424945745E-3", "-8.501026443778651E-3", "-8.2476738798851357E-3", "-6.3544557763970753E-3", "-6.6019744496516432E-3", "-6.2543887183258255E-3", "-6.5943053595595887E-3", "-8.4703470769107794E-3",
    "-4.6525507699802073E-3", "-7.7910337170988466E-3", "-4.5218055833669297E-3", "-6.4212381337125507E-3", "-6.2337616107754311E-3", "-8.208780527314578E-3", "-7.4015650604095226E-3", "-9.8147662065356869E-3", "-1.3826710276341388E-2", "-5.8696024390578741E-3",
    "-6.5859715869797525E-3", "-1.6486934187737248E-2", "-6.0692137372621441E-3", "-7.5571372309349927E-3", "-6.3493497516488463E-3", "-9.208860578361509E-3", "-1.3111641670367222E-2", "-1.1948918369000727E-2", "-8.5250657199378889E-3", "-7.3859973140404612E-3",
    "-1.0623901253811281E-2", "-4.9640353454574121E-3", "-1.0610296698846994E-2", "-5.73244649359754E-3", "-6.2197578471803365E-3", "-7.9659564693536512E-3", "-5.3758955235346274E-3", "-9.7173916344863668E-3", "-1.2317404329415095E-2", "-1.5127456752527305E-2"
)

for ($c = 0; $c -lt $row1Values.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = [double]$row1Values[$c]
    $ws.Cells.Item(2, $c + 1).Value = [double]$row2Values[$c]
}
